$d = $word.ActiveDocument

# Locate the "Git clone <URL of Repository>" list item - the new
# "Git branch <new branch name>" bullet belongs right after it.
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "*clone <URL of Repository>*") {
        $targetIndex = $i
        break
    }
}
if ($targetIndex -eq -1) {
    throw "Could not find the 'Git clone <URL of Repository>' paragraph"
}
$clonePara = $d.Paragraphs.Item($targetIndex)

# The document's trailing _GoBack bookmark currently sits at the end of the
# "Git clone" paragraph. Once we add the new bullet after it, that bookmark
# belongs at the end of the new (now last-edited) paragraph instead, so drop
# it here and recreate it in the right spot once the new text exists.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# Insert a new paragraph right after it; it inherits the ListParagraph /
# numbering formatting automatically from the paragraph mark it split off.
$clonePara.Range.InsertParagraphAfter() | Out-Null

$newPara = $d.Paragraphs.Item($targetIndex + 1)
$insertionRange = $newPara.Range.Duplicate
$insertionRange.MoveEnd(1, -1) | Out-Null   # exclude the paragraph mark

# Build the new bullet as real OOXML so it matches the run/proofErr
# structure Word itself produces ("Git" flagged & cleared by the spell
# checker, then the rest of the command text as a second run), with the
# _GoBack bookmark wrapping the very end of the new text.
$newParagraphXml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Git</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> branch &lt;new branch name&gt;</w:t></w:r><w:bookmarkStart w:id="100" w:name="_GoBack"/><w:bookmarkEnd w:id="100"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$insertionRange.InsertXML($newParagraphXml)
